# Updated cryptos list on Thu Apr  6 06:43:14 UTC 2023 with GitHub Actions
# Applies refreshed Price (column D) and Volume(1h) (column E) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.076.68"
$ws.Range("E2").Value = "'  -1.77%  "
$ws.Range("D3").Value = "'1.891.60"
$ws.Range("E4").Value = "'  -0.11%  "
$ws.Range("D5").Value = "'314.03"
$ws.Range("E5").Value = "'  -0.51%  "
$ws.Range("E6").Value = "'  -0.09%  "
$ws.Range("D7").Value = "'0.5048"
$ws.Range("E7").Value = "'  -2.26%  "
$ws.Range("D8").Value = "'0.3897"
$ws.Range("E8").Value = "'  -1.96%  "
$ws.Range("D9").Value = "'0.09245"
$ws.Range("E9").Value = "'  -4.88%  "
$ws.Range("E10").Value = "'  -2.63%  "
$ws.Range("D11").Value = "'41.81"
$ws.Range("E11").Value = "'  -0.33%  "
$ws.Range("D12").Value = "'6.386"
$ws.Range("E12").Value = "'  -2.47%  "
$ws.Range("D13").Value = "'20.82"
$ws.Range("E13").Value = "'  -2.15%  "
$ws.Range("D14").Value = "'1.891.59"
$ws.Range("E14").Value = "'  -1.24%  "
$ws.Range("E15").Value = "'  -4.03%  "
$ws.Range("E16").Value = "'  -0.09%  "
$ws.Range("D17").Value = "'92.19"
$ws.Range("E17").Value = "'  -1.95%  "
$ws.Range("E18").Value = "'  -2.89%  "
$ws.Range("D19").Value = "'0.06649"
$ws.Range("E19").Value = "'  -0.11%  "
$ws.Range("E20").Value = "'  -1.72%  "
$ws.Range("E21").Value = "'  +0.01%  "
$ws.Range("D22").Value = "'6.210"
$ws.Range("E22").Value = "'  -1.94%  "
$ws.Range("D23").Value = "'28.131.91"
$ws.Range("D24").Value = "'11.41"
$ws.Range("E24").Value = "'  -0.59%  "
$ws.Range("D25").Value = "'2.317"
$ws.Range("E25").Value = "'  +0.92%  "
$ws.Range("D26").Value = "'2.109.39"
$ws.Range("E26").Value = "'  -1.13%  "
$ws.Range("D27").Value = "'2.541"
$ws.Range("E27").Value = "'  -5.69%  "
$ws.Range("D28").Value = "'158.46"
$ws.Range("E28").Value = "'  -1.03%  "
$ws.Range("D29").Value = "'20.83"
$ws.Range("E29").Value = "'  -2.04%  "
$ws.Range("D30").Value = "'126.86"
$ws.Range("E30").Value = "'  -1.61%  "
$ws.Range("E31").Value = "'  -2.54%  "
$ws.Range("D32").Value = "'0.1054"
$ws.Range("E32").Value = "'  -2.75%  "
$ws.Range("E33").Value = "'  -2.71%  "
$ws.Range("D34").Value = "'3.605"
$ws.Range("E34").Value = "'  -1.00%  "
$ws.Range("D35").Value = "'9.464"
$ws.Range("E35").Value = "'  -4.47%  "
$ws.Range("D36").Value = "'0.06612"
$ws.Range("E36").Value = "'  -2.86%  "
$ws.Range("E37").Value = "'  +12.60%  "
$ws.Range("D38").Value = "'0.02403"
$ws.Range("D39").Value = "'0.2198"
$ws.Range("E39").Value = "'  -1.41%  "
$ws.Range("E40").Value = "'  -4.05%  "
$ws.Range("D41").Value = "'11.64"
$ws.Range("E41").Value = "'  -2.05%  "
$ws.Range("D42").Value = "'0.6441"
$ws.Range("E42").Value = "'  +0.00%  "
$ws.Range("D43").Value = "'4.970"
$ws.Range("E43").Value = "'  -3.05%  "
$ws.Range("D44").Value = "'0.9998"
$ws.Range("E44").Value = "'  -0.04%  "
$ws.Range("D45").Value = "'13.29"
$ws.Range("E45").Value = "'  -2.20%  "
$ws.Range("D46").Value = "'0.6053"
$ws.Range("E46").Value = "'  -0.61%  "
$ws.Range("E47").Value = "'  +1.42%  "
$ws.Range("D48").Value = "'3.687"
$ws.Range("E48").Value = "'  -2.48%  "
$ws.Range("D49").Value = "'2.004"
$ws.Range("E49").Value = "'  -1.78%  "
$ws.Range("D50").Value = "'122.05"
$ws.Range("E50").Value = "'  -2.59%  "
$ws.Range("E51").Value = "'  -1.64%  "
